$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attendance count
$ws.Range("A2").Value = "0/18"

# Update attendance status for students who are now absent
$ws.Range("B3").Value = "Fraværende"
$ws.Range("D3").Value = "Fraværende"
$ws.Range("M3").Value = "Fraværende"
